$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''29.501.50'
$ws.Range('E2').Value = '''  +2.09%  '
$ws.Range('D3').Value = '''1.857.36'
$ws.Range('E3').Value = '''  +1.29%  '
$ws.Range('D4').Value = '''0.9994'
$ws.Range('E4').Value = '''  +0.00%  '
$ws.Range('D5').Value = '''245.50'
$ws.Range('E5').Value = '''  +0.15%  '
$ws.Range('D6').Value = '''0.6961'
$ws.Range('E6').Value = '''  +1.11%  '
$ws.Range('D7').Value = '''1.000'
$ws.Range('E7').Value = '''  +0.04%  '
$ws.Range('E8').Value = '''  +0.72%  '
$ws.Range('D9').Value = '''0.07697'
$ws.Range('E9').Value = '''  +0.09%  '
$ws.Range('D10').Value = '''23.63'
$ws.Range('E10').Value = '''  +0.88%  '
$ws.Range('D11').Value = '''0.07791'
$ws.Range('E11').Value = '''  -0.36%  '
$ws.Range('D12').Value = '''5.165'
$ws.Range('E12').Value = '''  +1.32%  '
$ws.Range('D13').Value = '''1.851.77'
$ws.Range('E13').Value = '''  +1.11%  '
$ws.Range('D14').Value = '''0.6941'
$ws.Range('E14').Value = '''  +1.74%  '
$ws.Range('D15').Value = '''91.23'
$ws.Range('E15').Value = '''  +0.87%  '
$ws.Range('D16').Value = '''6.342'
$ws.Range('E16').Value = '''  -1.58%  '
$ws.Range('D17').Value = '''29.483.50'
$ws.Range('E17').Value = '''  +2.07%  '
$ws.Range('D18').Value = '''0.000008334'
$ws.Range('E18').Value = '''  +0.33%  '
$ws.Range('D19').Value = '''2.102.19'
$ws.Range('E19').Value = '''  +1.32%  '
$ws.Range('D20').Value = '''238.95'
$ws.Range('E20').Value = '''  -1.53%  '
$ws.Range('E21').Value = '''  +0.17%  '
$ws.Range('D22').Value = '''0.9994'
$ws.Range('E22').Value = '''  -0.06%  '
$ws.Range('D23').Value = '''7.623'
$ws.Range('E23').Value = '''  +2.09%  '
$ws.Range('D24').Value = '''1.000'
$ws.Range('E24').Value = '''  +0.06%  '
$ws.Range('E25').Value = '''  +1.41%  '
$ws.Range('D26').Value = '''160.05'
$ws.Range('E26').Value = '''  -0.82%  '
$ws.Range('D27').Value = '''8.893'
$ws.Range('E27').Value = '''  +0.87%  '
$ws.Range('D28').Value = '''18.27'
$ws.Range('E28').Value = '''  +0.37%  '
$ws.Range('D29').Value = '''1.532'
$ws.Range('E29').Value = '''  -0.88%  '
$ws.Range('D30').Value = '''4.251'
$ws.Range('E30').Value = '''  +0.82%  '
$ws.Range('D31').Value = '''4.151'
$ws.Range('E31').Value = '''  -0.12%  '
$ws.Range('D32').Value = '''1.204'
$ws.Range('E32').Value = '''  +1.63%  '
$ws.Range('D33').Value = '''0.05111'
$ws.Range('E33').Value = '''  -0.03%  '
$ws.Range('D34').Value = '''0.7774'
$ws.Range('E34').Value = '''  +1.52%  '
$ws.Range('D35').Value = '''1.881'
$ws.Range('E35').Value = '''  +1.97%  '
$ws.Range('D36').Value = '''1.149'
$ws.Range('E36').Value = '''  +0.65%  '
$ws.Range('D37').Value = '''2.688'
$ws.Range('E37').Value = '''  -0.06%  '
$ws.Range('D38').Value = '''1.316.62'
$ws.Range('E38').Value = '''  +7.51%  '
$ws.Range('D39').Value = '''0.01877'
$ws.Range('E39').Value = '''  +1.53%  '
$ws.Range('D40').Value = '''2.724'
$ws.Range('E40').Value = '''  +0.85%  '
$ws.Range('D41').Value = '''0.9537'
$ws.Range('E41').Value = '''  +1.27%  '
$ws.Range('D42').Value = '''106.09'
$ws.Range('E42').Value = '''  -2.41%  '
$ws.Range('D43').Value = '''5.774'
$ws.Range('E43').Value = '''  +1.07%  '
$ws.Range('D44').Value = '''1.001'
$ws.Range('E44').Value = '''  +0.14%  '
$ws.Range('D45').Value = '''9.842'
$ws.Range('E45').Value = '''  +3.15%  '
$ws.Range('B46').Value = 'RocketPoolETH'
$ws.Range('C46').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D46').Value = '''2.002.30'
$ws.Range('E46').Value = '''  +1.43%  '
$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D47').Value = '''0.5235'
$ws.Range('E47').Value = '''  +1.32%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').Value = '''1.788'
$ws.Range('E48').Value = '''  +2.31%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').Value = '''63.18'
$ws.Range('E49').Value = '''  -1.78%  '
$ws.Range('B50').Value = 'Aptos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D50').Value = '''6.961'
$ws.Range('E50').Value = '''  +0.88%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = '''0.05928'
$ws.Range('E51').Value = '''  +1.09%  '
